$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$subject = "Prueba en Excel "
$from = "Joan Martinez <joan_martinez.olivares@hotmail.com>"
$body = "Prueba en Excel 1 Obtener Outlook para iOS<https://aka.ms/o0ukef> "

for ($r = 4; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $subject
    $ws.Cells.Item($r, 2).Value = $from
    $ws.Cells.Item($r, 3).Value = $body
}
